$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. New column P needs to pick up the plain data-row border style (column O's
#    style) everywhere the table already has rows, before it gets values.
# ---------------------------------------------------------------------------
$ws.Range("O3").Copy()
$ws.Range("P3:P5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 1. Row 4 (2nd data row): O4 text changes from the "...\weights" string to
#    a plain "Training" result-path string, and a new P4 "Testing" result
#    path is added next to it.
# ---------------------------------------------------------------------------
$ws.Range("O4").Value = "runs\detect\train19"
$ws.Range("P4").Value = "runs\detect\train192"

# ---------------------------------------------------------------------------
# 2. Row 5 (3rd data row) gets filled in with real figures (it was almost
#    empty before).
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 817
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 77
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 0
$ws.Range("O5").Value = "runs\detect\train20"
$ws.Range("P5").Value = "runs\detect\train202"

# Match the formatting Excel gives the accuracy/precision figures once a row
# has usable data (numeric 0.0000 format instead of the plain/error style).
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("M5:N5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Two new rows (4th & 5th record) get appended below, copying row 3's
#    formatting (plain bordered style, no fill) and the same relative
#    formulas used throughout the table.
# ---------------------------------------------------------------------------
$ws.Range("A3:P3").Copy()
$ws.Range("A6:P6").PasteSpecial(-4122)
$ws.Range("A7:P7").PasteSpecial(-4122)

$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

$ws.Range("D6").Formula = "=ROUNDUP(C6*B6,0)"
$ws.Range("D7").Formula = "=ROUNDUP(C7*B7,0)"
$ws.Range("F6").Formula = "=B6-D6"
$ws.Range("F7").Formula = "=B7-D7"
$ws.Range("L6").Formula = "=H6+I6+J6+K6"
$ws.Range("L7").Formula = "=H7+I7+J7+K7"
$ws.Range("M6").Formula = "=((H6+I6)/(H6+I6+J6+K6))*100"
$ws.Range("M7").Formula = "=((H7+I7)/(H7+I7+J7+K7))*100"
$ws.Range("N6").Formula = "=(H6/(H6+I6))*100"
$ws.Range("N7").Formula = "=(H7/(H7+I7))*100"

# ---------------------------------------------------------------------------
# 4. "Result Saved" header (O1) now spans two columns (O1:P1) since there is
#    a Training-run path and a Testing-run path; O2/P2 carry the new
#    sub-header labels.
# ---------------------------------------------------------------------------
$ws.Range("O1:O2").UnMerge()
$ws.Range("O1:P1").Merge()
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

$ws.Range("O2").Value = "Training"
$ws.Range("P2").Value = "Testing"
$ws.Range("O1").Copy()
$ws.Range("O2:P2").PasteSpecial(-4122)
$ws.Range("O2:P2").VerticalAlignment = -4107
$ws.Range("O2").Value = "Training"
$ws.Range("P2").Value = "Testing"

# Column O keeps a similar (slightly narrower) width and the new column P
# mirrors it.
$ws.Columns.Item(15).ColumnWidth = 20.88
$ws.Columns.Item(16).ColumnWidth = 20.88

# ---------------------------------------------------------------------------
# 5. Selection moves to where the user was last working.
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
